$wb = $excel.ActiveWorkbook

# Add a new worksheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "7__c0904331-c00d-39b"

# Header row: copy the exact header formatting (bold, centered, top-aligned,
# thin border) used on every other sheet, then set the header text/values.
$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Range("A1:D1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("A1").Value = "rows"
$newSheet.Range("B1").Value = "columns"
$newSheet.Range("C1").Value = "year"
$newSheet.Range("D1").Value = "value"

# Data rows
$newSheet.Range("A2").Value = "None Cash at bank and in hand"
$newSheet.Range("C2").Value = 2023
$newSheet.Range("D2").Value = 14502614

$newSheet.Range("A3").Value = "None Term deposit"
$newSheet.Range("C3").Value = 2023
$newSheet.Range("D3").Value = 2368905

$newSheet.Range("A4").Value = "None Total"
$newSheet.Range("C4").Value = 2023
$newSheet.Range("D4").Value = 16871519

$newSheet.Range("A5").Value = "None Cash at bank and in hand"
$newSheet.Range("C5").Value = 2022
$newSheet.Range("D5").Value = 0

$newSheet.Range("A6").Value = "None Term deposit"
$newSheet.Range("C6").Value = 2022
$newSheet.Range("D6").Value = 0

$newSheet.Range("A7").Value = "None Total"
$newSheet.Range("C7").Value = 2022
$newSheet.Range("D7").Value = 13525659
